$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.896.03"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.690.59"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "'315.79"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "'0.3952"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("D8").Value = "'0.3985"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'1.443"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "'52.50"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").Value = "'1.011"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "'0.08721"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "'25.51"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "'7.374"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "'0.00001338"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "'7.858"
$ws.Range("E16").Value = "  -3.37%  "
$ws.Range("D17").Value = "1.712.47"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "'94.84"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").Value = "'0.07238"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "'20.39"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'7.149"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").Value = "'1.008"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'14.12"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "24.831.85"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "'2.390"
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("D26").Value = "'2.840"
$ws.Range("E26").Value = "  -4.92%  "
$ws.Range("D27").Value = "'23.17"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").Value = "'5.952"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "'161.49"
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").Value = "'148.77"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'8.081"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'2.624"
$ws.Range("E32").Value = "  +20.01%  "
$ws.Range("D33").Value = "1.898.73"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").Value = "'0.08480"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("D35").Value = "'0.03103"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "'1.029"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "'7.009"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").Value = "'0.2842"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").Value = "'0.09622"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("D40").Value = "'10.78"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "'0.8069"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'13.95"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'1.470"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "'16.89"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "'2.620"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "'0.7256"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'4.215"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'0.08930"
$ws.Range("E48").Value = "  +8.83%  "
$ws.Range("D49").Value = "'1.379"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'1.007"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "'138.94"
$ws.Range("E51").Value = "  -0.50%  "
